# Update "想去人数" (F column) counts that changed between the two
# generated data snapshots, on both the "展览" sheet and the "全部类型" sheet.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 140
$ws1.Range("F4").Value  = 2101
$ws1.Range("F5").Value  = 376
$ws1.Range("F6").Value  = 659
$ws1.Range("F8").Value  = 2087
$ws1.Range("F9").Value  = 10834
$ws1.Range("F12").Value = 293
$ws1.Range("F15").Value = 9054
$ws1.Range("F17").Value = 735
$ws1.Range("F18").Value = 5317
$ws1.Range("F20").Value = 3366
$ws1.Range("F21").Value = 3

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 140
$ws4.Range("F4").Value  = 2101
$ws4.Range("F5").Value  = 376
$ws4.Range("F6").Value  = 659
$ws4.Range("F9").Value  = 2087
$ws4.Range("F12").Value = 10834
$ws4.Range("F15").Value = 293
$ws4.Range("F18").Value = 9054
$ws4.Range("F20").Value = 735
$ws4.Range("F21").Value = 5317
$ws4.Range("F23").Value = 3366
$ws4.Range("F24").Value = 3
